$wb = $excel.ActiveWorkbook

$urlBase = "https://github.com/OpenLocalizationTestOrg/ol-test1/blob/4118363504852140ba86ce024d54b404c07eb606/e2e/"
$url9d = $urlBase + "9d6033cf-c2d7-4cbb-87ef-1e855a033654.md"
$urlD6 = $urlBase + "d6bed0a7-3bbd-49cf-8899-cc95e28d83c9.md"

# ---------------------------------------------------------------------------
# Sheet "Overview": rows for the two handback files swap position (row 2 now
# shows d6bed0a7's info, row 3 now shows 9d6033cf's info) and the 9d6033cf
# row picks up the new "Ready for handoff" status / timestamp.
# ---------------------------------------------------------------------------
$wsOverview = $wb.Worksheets.Item("Overview")

$wsOverview.Range("A2").Value = "d6bed0a7-3bbd-49cf-8899-cc95e28d83c9.md"
$wsOverview.Range("B2").Value = "e2e\d6bed0a7-3bbd-49cf-8899-cc95e28d83c9.md"

$wsOverview.Range("A3").Value = "9d6033cf-c2d7-4cbb-87ef-1e855a033654.md"
$wsOverview.Range("B3").Value = "e2e\9d6033cf-c2d7-4cbb-87ef-1e855a033654.md"

$wsOverview.Range("E3").Value = "Ready for handoff"
$wsOverview.Range("F3").Value = "Ready for handoff"
$wsOverview.Range("G3").Value = "2017-01-03 08:25:43"

$wsOverview.Columns.Item(5).ColumnWidth = 16.3
$wsOverview.Columns.Item(6).ColumnWidth = 16.3

# Hyperlinks: r:id keeps pointing at the same GitHub blob URL per row, but the
# visible display text now follows the swapped row content.
$wsOverview.Hyperlinks.Delete()
$wsOverview.Hyperlinks.Add($wsOverview.Range("B2"), $url9d, [Type]::Missing, [Type]::Missing, "e2e\d6bed0a7-3bbd-49cf-8899-cc95e28d83c9.md") | Out-Null
$wsOverview.Hyperlinks.Add($wsOverview.Range("B3"), $urlD6, [Type]::Missing, [Type]::Missing, "e2e\9d6033cf-c2d7-4cbb-87ef-1e855a033654.md") | Out-Null

# ---------------------------------------------------------------------------
# Sheet "zh-cn": same row swap, plus the 9d6033cf row's localization fields
# move from "In Translation"/ht to "Ready for handoff"/mt with new handoff
# timestamp.
# ---------------------------------------------------------------------------
$wsZhCn = $wb.Worksheets.Item("zh-cn")

$wsZhCn.Range("A2").Value = "d6bed0a7-3bbd-49cf-8899-cc95e28d83c9.md"
$wsZhCn.Range("B2").Value = ".md"
$wsZhCn.Range("C2").Value = "In Translation"
$wsZhCn.Range("D2").Value = "e2e"
$wsZhCn.Range("E2").Value = "ht"
$wsZhCn.Range("G2").Value = "d6bed0a7-3bbd-49cf-8899-cc95e28d83c9.0353b9c9c0654572a9e0e29ccb82bac604e33d9b.zh-cn.xlf"
$wsZhCn.Range("H2").Value = "2017-01-03 08:24:45"

$wsZhCn.Range("A3").Value = "9d6033cf-c2d7-4cbb-87ef-1e855a033654.md"
$wsZhCn.Range("B3").Value = ".md"
$wsZhCn.Range("C3").Value = "Ready for handoff"
$wsZhCn.Range("D3").Value = "e2e"
$wsZhCn.Range("E3").Value = "mt"
$wsZhCn.Range("G3").Value = "9d6033cf-c2d7-4cbb-87ef-1e855a033654.aa3fee4501fb3f951ec6bcf3ef6188d8accb02d2.zh-cn.xlf"
$wsZhCn.Range("H3").Value = "2017-01-03 08:25:34"

$wsZhCn.Columns.Item(3).ColumnWidth = 16.3

$wsZhCn.Hyperlinks.Delete()
$wsZhCn.Hyperlinks.Add($wsZhCn.Range("A2"), $url9d, [Type]::Missing, [Type]::Missing, "d6bed0a7-3bbd-49cf-8899-cc95e28d83c9.md") | Out-Null
$wsZhCn.Hyperlinks.Add($wsZhCn.Range("A3"), $urlD6, [Type]::Missing, [Type]::Missing, "9d6033cf-c2d7-4cbb-87ef-1e855a033654.md") | Out-Null

# ---------------------------------------------------------------------------
# Sheet "de-de": same row swap; the 9d6033cf row status flips too, but
# Priority stays "ht" and only the handoff timestamp is refreshed.
# ---------------------------------------------------------------------------
$wsDeDe = $wb.Worksheets.Item("de-de")

$wsDeDe.Range("A2").Value = "d6bed0a7-3bbd-49cf-8899-cc95e28d83c9.md"
$wsDeDe.Range("B2").Value = ".md"
$wsDeDe.Range("C2").Value = "In Translation"
$wsDeDe.Range("D2").Value = "e2e"
$wsDeDe.Range("E2").Value = "ht"
$wsDeDe.Range("G2").Value = "d6bed0a7-3bbd-49cf-8899-cc95e28d83c9.0353b9c9c0654572a9e0e29ccb82bac604e33d9b.de-de.xlf"
$wsDeDe.Range("H2").Value = "2017-01-03 08:24:54"

$wsDeDe.Range("A3").Value = "9d6033cf-c2d7-4cbb-87ef-1e855a033654.md"
$wsDeDe.Range("B3").Value = ".md"
$wsDeDe.Range("C3").Value = "Ready for handoff"
$wsDeDe.Range("D3").Value = "e2e"
$wsDeDe.Range("E3").Value = "ht"
$wsDeDe.Range("G3").Value = "9d6033cf-c2d7-4cbb-87ef-1e855a033654.aa3fee4501fb3f951ec6bcf3ef6188d8accb02d2.de-de.xlf"
$wsDeDe.Range("H3").Value = "2017-01-03 08:25:43"

$wsDeDe.Columns.Item(3).ColumnWidth = 16.3

$wsDeDe.Hyperlinks.Delete()
$wsDeDe.Hyperlinks.Add($wsDeDe.Range("A2"), $url9d, [Type]::Missing, [Type]::Missing, "d6bed0a7-3bbd-49cf-8899-cc95e28d83c9.md") | Out-Null
$wsDeDe.Hyperlinks.Add($wsDeDe.Range("A3"), $urlD6, [Type]::Missing, [Type]::Missing, "9d6033cf-c2d7-4cbb-87ef-1e855a033654.md") | Out-Null
